$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new rows before row 231, shifting existing rows 231-233 down to 233-235
$ws.Rows.Item(231).Insert()
$ws.Rows.Item(231).Insert()

# Row 231: new match Zalaegerszegi TE vs Kecskemeti TE
$ws.Cells.Item(231,1).Value = 229.0
$ws.Cells.Item(231,2).Value = 6818327.0
$ws.Cells.Item(231,3).Value = "Hungary NB I"
$ws.Cells.Item(231,4).Value = "Hungary NB I"
$ws.Cells.Item(231,5).Value = 45346.35416666666
$ws.Cells.Item(231,6).Value = "Zalaegerszegi TE"
$ws.Cells.Item(231,7).Value = "Kecskemeti TE"
$ws.Cells.Item(231,8).Value = 3.0
$ws.Cells.Item(231,9).Value = 1.0
$ws.Cells.Item(231,10).Value = "H"
$ws.Cells.Item(231,11).Value = 2.4
$ws.Cells.Item(231,12).Value = 3.2
$ws.Cells.Item(231,13).Value = 2.8
$ws.Cells.Item(231,14).Value = 2.375
$ws.Cells.Item(231,15).Value = 3.0
$ws.Cells.Item(231,16).Value = 3.1
$ws.Cells.Item(231,17).Value = -0.25
$ws.Cells.Item(231,18).Value = 2.05
$ws.Cells.Item(231,19).Value = 1.8
$ws.Cells.Item(231,20).Value = 2.25
$ws.Cells.Item(231,21).Value = 1.8
$ws.Cells.Item(231,22).Value = 2.05
$ws.Cells.Item(231,23).Value = 1.375
$ws.Cells.Item(231,24).Value = -1.0
$ws.Cells.Item(231,25).Value = -1.0
$ws.Cells.Item(231,26).Value = 1.05
$ws.Cells.Item(231,27).Value = -1.0
$ws.Cells.Item(231,28).Value = 0.8
$ws.Cells.Item(231,29).Value = -1.0

# Row 232: new match Diosgyori VTK vs MOL Fehervar FC
$ws.Cells.Item(232,1).Value = 230.0
$ws.Cells.Item(232,2).Value = 6818326.0
$ws.Cells.Item(232,3).Value = "Hungary NB I"
$ws.Cells.Item(232,4).Value = "Hungary NB I"
$ws.Cells.Item(232,5).Value = 45346.45833333334
$ws.Cells.Item(232,6).Value = "Diosgyori VTK"
$ws.Cells.Item(232,7).Value = "MOL Fehervar FC"
$ws.Cells.Item(232,8).Value = 4.0
$ws.Cells.Item(232,9).Value = 0.0
$ws.Cells.Item(232,10).Value = "H"
$ws.Cells.Item(232,11).Value = 2.25
$ws.Cells.Item(232,12).Value = 3.6
$ws.Cells.Item(232,13).Value = 2.75
$ws.Cells.Item(232,14).Value = 1.95
$ws.Cells.Item(232,15).Value = 3.6
$ws.Cells.Item(232,16).Value = 3.5
$ws.Cells.Item(232,17).Value = -0.5
$ws.Cells.Item(232,18).Value = 2.05
$ws.Cells.Item(232,19).Value = 1.8
$ws.Cells.Item(232,20).Value = 2.75
$ws.Cells.Item(232,21).Value = 1.975
$ws.Cells.Item(232,22).Value = 1.875
$ws.Cells.Item(232,23).Value = 0.95
$ws.Cells.Item(232,24).Value = -1.0
$ws.Cells.Item(232,25).Value = -1.0
$ws.Cells.Item(232,26).Value = 1.05
$ws.Cells.Item(232,27).Value = -1.0
$ws.Cells.Item(232,28).Value = 0.9750000000000001
$ws.Cells.Item(232,29).Value = -1.0

# Row 233 (former row 231, now with result filled in)
$ws.Cells.Item(233,8).Value = 1.0
$ws.Cells.Item(233,9).Value = 0.0
$ws.Cells.Item(233,10).Value = "H"
$ws.Cells.Item(233,23).Value = 1.45
$ws.Cells.Item(233,24).Value = -1.0
$ws.Cells.Item(233,25).Value = -1.0
$ws.Cells.Item(233,26).Value = 0.825
$ws.Cells.Item(233,27).Value = -1.0
$ws.Cells.Item(233,28).Value = -1.0
$ws.Cells.Item(233,29).Value = 1.0

# Row 234 (former row 232, odds refreshed)
$ws.Cells.Item(234,18).Value = 1.85
$ws.Cells.Item(234,19).Value = 2.0
$ws.Cells.Item(234,21).Value = 1.85
$ws.Cells.Item(234,22).Value = 2.0

